# Hortaliza, Feria Lagunitas de Puerto Montt - Perejil
# A new weekly price record is inserted as row 82, pushing the existing
# rows 82..191 down to 83..192 (dimension grows from A1:R191 to A1:R192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 82, shifting rows 82-191
# down to 83-192 (row 191's old data ends up duplicated into new row 192).
$ws.Range("A82:R82").Insert()

# Populate the freshly inserted row 82 with the new record.
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44546
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112044
$ws.Range("G82").Value = "Perejil"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 30
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = 5000
$ws.Range("N82").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 1667
$ws.Range("Q82").Value = 3
$ws.Range("R82").Value = "Hortaliza"
